$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3448735177516937
$ws.Range("C2").Value = 0.3289481997489929

$ws.Range("B3").Value = 0.4984879726753482
$ws.Range("C3").Value = 0.5159719049806877

$ws.Range("B4").Value = 0.2227258788774616
$ws.Range("C4").Value = 0.2120054008814109

$ws.Range("B5").Value = 0.2230000048875809
$ws.Range("C5").Value = 0.1959999948740005

$ws.Range("B6").Value = -0.2570435404777527
$ws.Range("C6").Value = -0.278084397315979
